# Applies the "Updated symbol list" data refresh (Tue Jan 31 23:38:30 UTC 2023)
# to the crypto-tracker worksheet: each row's Price (D) / Volume 1h (E) figures
# move to their newly-scraped values, and several rows shift which coin/link
# (B/C) they describe because the source ranking re-sorted.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Every cell in this sheet was authored as literal text (coinranking.com's
# exporter writes prices/percentages as strings, not numbers), so re-apply the
# same convention: force Text number format before the assignment (otherwise
# COM auto-converts "311.61" / "1.77%" into a real number/percentage and we'd
# lose the exact original formatting such as trailing zeros or the % sign),
# then drop the formatting back to Normal so no stray style is left behind.
$updates = @(
    @{ Cell = "D2"; Value = "311.61" },
    @{ Cell = "E2"; Value = "1.77%" },
    @{ Cell = "E3"; Value = "0.71%" },
    @{ Cell = "D4"; Value = "5.122" },
    @{ Cell = "E4"; Value = "0.97%" },
    @{ Cell = "D5"; Value = "0.07876" },
    @{ Cell = "E5"; Value = "1.87%" },
    @{ Cell = "B6"; Value = "FTXToken" },
    @{ Cell = "C6"; Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt" },
    @{ Cell = "D6"; Value = "1.905" },
    @{ Cell = "E6"; Value = "0.22%" },
    @{ Cell = "D7"; Value = "8.266" },
    @{ Cell = "B8"; Value = "BTSEToken" },
    @{ Cell = "C8"; Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse" },
    @{ Cell = "D8"; Value = "2.838" },
    @{ Cell = "E8"; Value = "-7.50%" },
    @{ Cell = "B9"; Value = "MXToken" },
    @{ Cell = "C9"; Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx" },
    @{ Cell = "D9"; Value = "0.9203" },
    @{ Cell = "E9"; Value = "0.11%" },
    @{ Cell = "B10"; Value = "LiechtensteinCryptoassetsExchange" },
    @{ Cell = "C10"; Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx" },
    @{ Cell = "D10"; Value = "0.1173" },
    @{ Cell = "E10"; Value = "-4.40%" },
    @{ Cell = "B11"; Value = "WazirX" },
    @{ Cell = "C11"; Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx" },
    @{ Cell = "D11"; Value = "0.1931" },
    @{ Cell = "E11"; Value = "3.31%" },
    @{ Cell = "B12"; Value = "MandalaExchangeToken" },
    @{ Cell = "C12"; Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx" },
    @{ Cell = "D12"; Value = "0.09028" },
    @{ Cell = "E12"; Value = "2.51%" },
    @{ Cell = "B13"; Value = "BitrueCoin" },
    @{ Cell = "C13"; Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr" },
    @{ Cell = "D13"; Value = "0.03317" },
    @{ Cell = "E13"; Value = "-2.55%" },
    @{ Cell = "B14"; Value = "BitMartToken" },
    @{ Cell = "C14"; Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx" },
    @{ Cell = "D14"; Value = "0.09599" },
    @{ Cell = "E14"; Value = "-0.99%" },
    @{ Cell = "B15"; Value = "BitForexToken" },
    @{ Cell = "C15"; Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf" },
    @{ Cell = "D15"; Value = "0.001376" },
    @{ Cell = "E15"; Value = "0.90%" },
    @{ Cell = "B16"; Value = "TigerCash" },
    @{ Cell = "C16"; Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch" },
    @{ Cell = "D16"; Value = "0.006091" },
    @{ Cell = "E16"; Value = "2.22%" },
    @{ Cell = "B17"; Value = "LEO" },
    @{ Cell = "C17"; Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo" },
    @{ Cell = "D17"; Value = "3.539" },
    @{ Cell = "E17"; Value = "-0.92%" },
    @{ Cell = "B18"; Value = "GateToken" },
    @{ Cell = "C18"; Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt" },
    @{ Cell = "D18"; Value = "4.411" },
    @{ Cell = "E18"; Value = "1.65%" },
    @{ Cell = "D20"; Value = "5.263" },
    @{ Cell = "E20"; Value = "4.94%" },
    @{ Cell = "D21"; Value = "0.1284" },
    @{ Cell = "E21"; Value = "0.46%" },
    @{ Cell = "E22"; Value = "3.98%" },
    @{ Cell = "D23"; Value = "0.04352" },
    @{ Cell = "E23"; Value = "0.53%" },
    @{ Cell = "E24"; Value = "3.03%" },
    @{ Cell = "D25"; Value = "0.004665" },
    @{ Cell = "E25"; Value = "10.35%" },
    @{ Cell = "E26"; Value = "0.73%" },
    @{ Cell = "D27"; Value = "0.0003986" },
    @{ Cell = "D39"; Value = "0.02250" },
    @{ Cell = "E39"; Value = "3.68%" },
    @{ Cell = "D40"; Value = "0.05085" },
    @{ Cell = "E40"; Value = "3.70%" },
    @{ Cell = "D41"; Value = "0.007459" },
    @{ Cell = "E41"; Value = "-1.91%" },
    @{ Cell = "D42"; Value = "0.009027" },
    @{ Cell = "E42"; Value = "-9.46%" },
    @{ Cell = "D43"; Value = "0.1352" },
    @{ Cell = "E43"; Value = "0.94%" },
    @{ Cell = "D44"; Value = "0.001949" },
    @{ Cell = "E44"; Value = "-5.36%" },
    @{ Cell = "D45"; Value = "0.008627" },
    @{ Cell = "E45"; Value = "-12.23%" },
    @{ Cell = "D46"; Value = "0.00006554" },
    @{ Cell = "E46"; Value = "0.38%" },
    @{ Cell = "E47"; Value = "-0.02%" },
    @{ Cell = "B48"; Value = "CoinbaseStockToken" },
    @{ Cell = "C48"; Value = "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin" },
    @{ Cell = "D48"; Value = "0.0009991" },
    @{ Cell = "E48"; Value = "-23.19%" },
    @{ Cell = "B49"; Value = "BOLO" },
    @{ Cell = "C49"; Value = "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo" },
    @{ Cell = "D49"; Value = "0.003160" },
    @{ Cell = "E49"; Value = "5.51%" },
    @{ Cell = "E50"; Value = "-0.02%" },
    @{ Cell = "E51"; Value = "-0.02%" }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    $col = $u.Cell -replace '[0-9]+$', ''
    if ($col -eq "D" -or $col -eq "E") {
        $cell.NumberFormat = "@"
        $cell.Value = $u.Value
        $cell.Style = "Normal"
    } else {
        $cell.Value = $u.Value
    }
}
